$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Helper: replace the first occurrence of $oldVal inside a specific table
# cell with $newVal, by computing an absolute document Range so the edit
# stays confined to that one cell (Range.Find.Execute in this host ignores
# the owning range's boundaries and searches/replaces document-wide).
function Replace-CellText($table, $row, $col, $oldVal, $newVal) {
    $cell = $table.Cell($row, $col)
    $cellRange = $cell.Range
    $cellText = $cellRange.Text
    $pos = $cellText.IndexOf($oldVal)
    if ($pos -ge 0) {
        $target = $d.Range($cellRange.Start + $pos, $cellRange.Start + $pos + $oldVal.Length)
        $target.Text = $newVal
        return $true
    }
    return $false
}

# 1. CNPJ: formatted value -> digits-only value (single occurrence in doc)
$cnpjFound = $d.Content.Find.Execute("11.064.624/0001-99", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1106462000199", 2)
if (-not $cnpjFound) { Write-Host "WARNING: CNPJ text not found" }

# 2. Row 7, col 2 - "II - Revenda de mercadorias com documento fiscal emitido" -> R$ 3000.00
if (-not (Replace-CellText $t 7 2 "0.00" "3000.00")) { Write-Host "WARNING: row 7 value not found" }

# 3. Row 8, col 2 - "III - Total das receitas com revenda de mercadorias (I + II)" -> R$ 3000.00
if (-not (Replace-CellText $t 8 2 "0.00" "3000.00")) { Write-Host "WARNING: row 8 value not found" }

# 4. Row 17, col 2 - "X - Total geral das receitas brutas no mes (III + VI + IX)" -> R$ 3000.00
if (-not (Replace-CellText $t 17 2 "0.00" "3000.00")) { Write-Host "WARNING: row 17 value not found" }
